$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Append a new row (51) below the last existing data row (50), re-using its
# formatting (normal text style for most cells, hyperlink style for F/G).
$srcRange = $ws.Range("A50:J50")
$destRange = $ws.Range("A51:J51")
$srcRange.Copy($destRange)

$values = @(
    "Q_UBALAWA",
    "Umweltbundesamt nach Angaben der Bund/Länder Arbeitsgemeinschaft Wasser (LAWA)",
    "German Environment Agency on the basis of data from the German Working Group on Water Issues of the Länder and the Federal Government (LAWA)",
    "Umweltbundesamt nach Angaben der Bund/Länder Arbeitsgemeinschaft Wasser (LAWA)",
    "German Environment Agency on the basis of data from the German Working Group on Water Issues of the Länder and the Federal Government (LAWA)",
    "https://www.umweltbundesamt.de/",
    "https://www.umweltbundesamt.de/en",
    "",
    "",
    ""
)

for ($col = 1; $col -le 10; $col++) {
    $ws.Cells.Item(51, $col).Value = $values[$col - 1]
}
